$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.268.30'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '2.593.71'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '569.45'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').Value = '141.61'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '2.615.32'
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').Value = '6.58'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').Value = '0.369'
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('E13').Value = '  -6.37%  '
$ws.Range('D14').Value = '3.072.22'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '60.306.12'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '23.19'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').Value = '2.604.29'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = '11.35'
$ws.Range('E19').Value = '  +9.68%  '
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '346.25'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').Value = '6.98'
$ws.Range('E22').Value = '  +9.15%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').Value = '0.534'
$ws.Range('E24').Value = '  +14.64%  '
$ws.Range('D25').Value = '63.13'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('D28').Value = '7.68'
$ws.Range('E28').Value = '  +5.17%  '
$ws.Range('D29').Value = '0.0₃0784'
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('D30').Value = '1.83'
$ws.Range('E30').Value = '  +9.40%  '
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '6.34'
$ws.Range('E32').Value = '  +4.06%  '
$ws.Range('D33').Value = '160.91'
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E34').Value = '  +2.64%  '
$ws.Range('E35').Value = '  +5.06%  '
$ws.Range('E36').Value = '  +10.24%  '
$ws.Range('E37').Value = '  +4.64%  '
$ws.Range('D38').Value = '1.61'
$ws.Range('E38').Value = '  +8.89%  '
$ws.Range('D39').Value = '37.82'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.81'
$ws.Range('E40').Value = '  +4.05%  '
$ws.Range('B41').Value = 'SuiNetwork'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D41').Value = '0.851'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('D42').Value = '294.33'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').Value = '138.18'
$ws.Range('E43').Value = '  +4.94%  '
$ws.Range('D44').Value = '0.997'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').Value = '0.607'
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('D46').Value = '0.0979'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').Value = '19.69'
$ws.Range('E47').Value = '  +3.63%  '
$ws.Range('E48').Value = '  +1.95%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '19.80'
$ws.Range('E49').Value = '  +6.62%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0239'
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '10.67'
$ws.Range('E51').Value = '  +0.24%  '
